$d = $word.ActiveDocument

# Locate the existing " Gupta" run (it sits right before the _GoBack bookmark).
$matchRange = $d.Content
$found = $matchRange.Find.Execute(" Gupta", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# The matched range covers " Gupta" (a leading space + "Gupta").
# We intentionally keep the leading space character untouched (as the
# original run) and only replace "Gupta" (i.e. Range.Start + 1 .. Range.End).
# This avoids landing the insertion point exactly on top of the zero-length
# "_GoBack" bookmark that immediately follows this text, which would
# otherwise cause the bookmark to be split/expanded around the new content.
$target = $d.Range($matchRange.Start + 1, $matchRange.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body>' + `
            '<w:p>' + `
              '<w:proofErr w:type="spellStart"/>' + `
              '<w:r><w:t>Gupta</w:t></w:r>' + `
              '<w:r><w:t>.I</w:t></w:r>' + `
              '<w:proofErr w:type="spellEnd"/>' + `
              '<w:r><w:t xml:space="preserve"> am a student.</w:t></w:r>' + `
            '</w:p>' + `
          '</w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

$target.InsertXML($xml)
